# Apply "Added recursive and iterative sorts comparison" edit.
#
# Summary of changes:
#  - same_elements (sheet5) and partly_same (sheet6) gain a new "G" column
#    (header 500000) with freshly measured byte/int sort timings, so the
#    used range grows from A1:F3 to A1:G3 on both sheets and the existing
#    B:F timing values are refreshed to the newly measured numbers.
#  - Each sheet's view now remembers a cell/range selection.
#  - The workbook window remembers that sheet tab scrolling starts at the
#    second tab (firstSheet).

$wb = $excel.ActiveWorkbook

# Keep track of whichever sheet is active right now so we can restore it -
# selecting on a sheet implicitly activates it, and we do not want to change
# which tab is active overall.
$originalActive = $wb.ActiveSheet

$straight      = $wb.Worksheets.Item(1)   # straight
$reversed      = $wb.Worksheets.Item(2)   # reversed
$sorted        = $wb.Worksheets.Item(3)   # sorted
$partlySorted  = $wb.Worksheets.Item(4)   # partly_sorted
$sameElements  = $wb.Worksheets.Item(5)   # same_elements
$partlySame    = $wb.Worksheets.Item(6)   # partly_same

# ---------------------------------------------------------------------
# same_elements: add column G and refresh the measured timings
# ---------------------------------------------------------------------
$sameElements.Range("G1").Value = 500000

$sameElements.Range("B2").Value = 0.003503
$sameElements.Range("C2").Value = 0.01852
$sameElements.Range("D2").Value = 0.190199
$sameElements.Range("E2").Value = 1.871455
$sameElements.Range("F2").Value = 21.047983
$sameElements.Range("G2").Value = 206.405991

$sameElements.Range("B3").Value = 0.004505
$sameElements.Range("C3").Value = 0.017517
$sameElements.Range("D3").Value = 0.225736
$sameElements.Range("E3").Value = 1.839922
$sameElements.Range("F3").Value = 22.778308
$sameElements.Range("G3").Value = 229.372907

# ---------------------------------------------------------------------
# partly_same: add column G and refresh the measured timings
# ---------------------------------------------------------------------
$partlySame.Range("G1").Value = 500000

$partlySame.Range("B2").Value = 0.003003
$partlySame.Range("C2").Value = 0.018519
$partlySame.Range("D2").Value = 0.229739
$partlySame.Range("E2").Value = 2.420529
$partlySame.Range("F2").Value = 21.456925
$partlySame.Range("G2").Value = 221.834946

$partlySame.Range("B3").Value = 0.004505
$partlySame.Range("C3").Value = 0.022023
$partlySame.Range("D3").Value = 0.24926
$partlySame.Range("E3").Value = 2.188285
$partlySame.Range("F3").Value = 22.8789
$partlySame.Range("G3").Value = 229.676239

# ---------------------------------------------------------------------
# Sheet view selections - select the full used range A1:G3 on each of the
# first five sheets (mirrors the reviewer selecting the data after the
# new column landed).
# ---------------------------------------------------------------------
$straight.Activate()
$straight.Range("A1:G3").Select()

$reversed.Activate()
$reversed.Range("A1:G3").Select()

$sorted.Activate()
$sorted.Range("A1:G3").Select()

$partlySorted.Activate()
$partlySorted.Range("A1:G3").Select()

$sameElements.Activate()
$sameElements.Range("A1:G3").Select()

# partly_same keeps a single-cell selection, just moved to J23.
$partlySame.Activate()
$partlySame.Range("J23").Select()

# Restore the sheet that was active before we started touching selections.
$originalActive.Activate()

# ---------------------------------------------------------------------
# The window remembers that the tab strip is scrolled so the first
# visible sheet tab is the second one (firstSheet="1").
# ---------------------------------------------------------------------
try {
    $wb.Windows.Item(1).ScrollWorkbookTabs(1)
} catch {
}
